# Stock.xlsx — update Quantity column values:
#   C2 (PC quantity):      10 -> 4
#   C3 (Tomate quantity): 100 -> 95
#
# The source values are stored as text (shared strings) even though they
# look numeric, so we briefly mark the cells as Text before writing the new
# value (otherwise Excel auto-converts a numeric-looking entry to a real
# number) and then clear the formatting again so the cell keeps using the
# workbook's default (General) style, just like every other cell here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "4"
$ws.Range("C2").ClearFormats()

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "95"
$ws.Range("C3").ClearFormats()
